$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that sometimes look like plain numbers
# (e.g. "578.15") and sometimes don't (e.g. "65.111.64"). Excel's COM
# layer auto-detects numeric-looking text and stores it as a real number
# when assigned via .Value. Forcing the whole price column to Text format
# first keeps every write a literal string (matching the source data,
# which is inline/shared text, not numeric), then we restore the default
# style so no stray formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.111.64"
$ws.Range("E2").Value = "  +2.67%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.477.83"
$ws.Range("E3").Value = "  +2.56%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.17%  "

# Row 5 - BNB
$ws.Range("D5").Value = "578.15"
$ws.Range("E5").Value = "  +0.28%  "

# Row 6 - Solana
$ws.Range("D6").Value = "162.33"
$ws.Range("E6").Value = "  +4.65%  "

# Row 7 - USDC
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.482.12"
$ws.Range("E8").Value = "  +2.09%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").Value = "  +9.13%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "7.35"
$ws.Range("E10").Value = "  -2.26%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +4.36%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  +1.66%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.058.78"
$ws.Range("E13").Value = "  +1.98%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -2.91%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +5.72%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "28.87"
$ws.Range("E16").Value = "  +6.57%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "65.059.24"
$ws.Range("E17").Value = "  +2.46%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.442.09"
$ws.Range("E18").Value = "  +0.74%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  +0.29%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "14.36"
$ws.Range("E20").Value = "  +2.28%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "390.84"
$ws.Range("E21").Value = "  +1.00%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "8.23"
$ws.Range("E22").Value = "  -1.90%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.549"
$ws.Range("E23").Value = "  +2.70%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +1.07%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "73.07"
$ws.Range("E25").Value = "  +3.03%  "

# Row 26 - PEPE
$ws.Range("D26").Value = "0.0000125"
$ws.Range("E26").Value = "  +19.88%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "9.54"
$ws.Range("E27").Value = "  +0.01%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +0.48%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.02%  "

# Row 30 - NEARProtocol
$ws.Range("E30").Value = "  +9.00%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  +8.33%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +0.45%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "23.68"
$ws.Range("E33").Value = "  +2.50%  "

# Row 34 - RenderToken
$ws.Range("D34").Value = "6.53"
$ws.Range("E34").Value = "  -0.08%  "

# Row 36 - Aptos
$ws.Range("E36").Value = "  +6.00%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +1.35%  "

# Row 38 - Monero
$ws.Range("D38").Value = "161.70"
$ws.Range("E38").Value = "  +2.48%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +1.80%  "

# Row 40 - Maker
$ws.Range("D40").Value = "2.993.34"
$ws.Range("E40").Value = "  +2.10%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "0.0769"
$ws.Range("E41").Value = "  +0.17%  "

# Row 42 - EnergySwap
$ws.Range("D42").Value = "27.51"
$ws.Range("E42").Value = "  -0.17%  "

# Row 43 - Filecoin
$ws.Range("E43").Value = "  +6.35%  "

# Row 44 - OKB
$ws.Range("D44").Value = "42.91"
$ws.Range("E44").Value = "  +3.75%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  -1.20%  "

# Row 46 - Mantle
$ws.Range("D46").Value = "0.781"
$ws.Range("E46").Value = "  +1.86%  "

# Row 47 - InjectiveProtocol
$ws.Range("D47").Value = "24.27"
$ws.Range("E47").Value = "  +8.00%  "

# Row 48 - ONDO
$ws.Range("D48").Value = "1.10"
$ws.Range("E48").Value = "  +2.60%  "

# Row 49 - was dogwifhat, now SuiNetwork
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").Value = "0.877"
$ws.Range("E49").Value = "  +7.51%  "

# Row 50 - was SuiNetwork, now dogwifhat
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.20"
$ws.Range("E50").Value = "  +12.34%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +4.13%  "

# Restore the default (unstyled) formatting on the price column now that
# every write has happened, so no residual "Text" number format lingers.
$ws.Range("D2:D51").Style = "Normal"
